# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect refreshed counts captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1373
$wsExhibit.Range("F7").Value = 11737
$wsExhibit.Range("F8").Value = 4400
$wsExhibit.Range("F15").Value = 152
$wsExhibit.Range("F17").Value = 5116
$wsExhibit.Range("F21").Value = 11353
$wsExhibit.Range("F22").Value = 11293

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1373
$wsAll.Range("F7").Value = 11737
$wsAll.Range("F8").Value = 4400
$wsAll.Range("F16").Value = 152
$wsAll.Range("F18").Value = 5116
$wsAll.Range("F22").Value = 11353
$wsAll.Range("F23").Value = 11293
